# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh updates to the Diabolos_Profits workbook
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 969.0227
$ws.Cells.Item(15, 9).Value = 969.0227
$ws.Cells.Item(15, 11).Value = 2907.0681
$ws.Cells.Item(15, 13).Value = -2738.0681
$ws.Cells.Item(17, 8).Value = 387846.78
$ws.Cells.Item(17, 10).Value = 401760.62
$ws.Cells.Item(17, 12).Value = 1205281.86
$ws.Cells.Item(17, 14).Value = -1205617.86
$ws.Cells.Item(28, 8).Value = 46890.816
$ws.Cells.Item(28, 9).Value = 51300.4
$ws.Cells.Item(28, 10).Value = 2795
$ws.Cells.Item(28, 11).Value = 51300.4
$ws.Cells.Item(28, 12).Value = 2795
$ws.Cells.Item(28, 13).Value = -50815.4
$ws.Cells.Item(28, 14).Value = -3765
$ws.Cells.Item(33, 8).Value = 62928.082
$ws.Cells.Item(33, 9).Value = 94072
$ws.Cells.Item(33, 10).Value = 640.25
$ws.Cells.Item(33, 11).Value = 94072
$ws.Cells.Item(33, 12).Value = 640.25
$ws.Cells.Item(33, 13).Value = -93843
$ws.Cells.Item(33, 14).Value = -1098.25
$ws.Cells.Item(64, 8).Value = 5439.2
$ws.Cells.Item(64, 10).Value = 5855.25
$ws.Cells.Item(64, 12).Value = 5855.25
$ws.Cells.Item(64, 14).Value = -6351.25
$ws.Cells.Item(67, 8).Value = 5439.2
$ws.Cells.Item(67, 10).Value = 5855.25
$ws.Cells.Item(67, 12).Value = 5855.25
$ws.Cells.Item(67, 14).Value = -7571.25
$ws.Cells.Item(92, 8).Value = 41805.082
$ws.Cells.Item(92, 9).Value = 144.43478
$ws.Cells.Item(92, 11).Value = 144.43478
$ws.Cells.Item(92, 13).Value = 1103.56522
$ws.Cells.Item(137, 8).Value = 4003.5
$ws.Cells.Item(137, 9).Value = 3105
$ws.Cells.Item(137, 10).Value = 6159.9
$ws.Cells.Item(137, 11).Value = 9315
$ws.Cells.Item(137, 12).Value = 18479.7
$ws.Cells.Item(137, 13).Value = -6765
$ws.Cells.Item(137, 14).Value = -23579.7
$ws.Cells.Item(138, 8).Value = 3549.862
$ws.Cells.Item(138, 10).Value = 4042.5
$ws.Cells.Item(138, 12).Value = 12127.5
$ws.Cells.Item(138, 14).Value = -22407.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 22115.46
$ws.Cells.Item(32, 9).Value = 13332.512
$ws.Cells.Item(32, 11).Value = 13332.512
$ws.Cells.Item(32, 13).Value = -13045.512
$ws.Cells.Item(40, 8).Value = 18749
$ws.Cells.Item(40, 10).Value = 17665.666
$ws.Cells.Item(40, 12).Value = 17665.666
$ws.Cells.Item(40, 14).Value = -18017.666
$ws.Cells.Item(80, 8).Value = 96666.664
$ws.Cells.Item(80, 9).Value = 70000
$ws.Cells.Item(80, 11).Value = 70000
$ws.Cells.Item(80, 13).Value = -69002
$ws.Cells.Item(83, 8).Value = 96666.664
$ws.Cells.Item(83, 9).Value = 70000
$ws.Cells.Item(83, 11).Value = 210000
$ws.Cells.Item(83, 13).Value = -205008
$ws.Cells.Item(97, 8).Value = 900.35
$ws.Cells.Item(97, 9).Value = 842.4737
$ws.Cells.Item(97, 11).Value = 842.4737
$ws.Cells.Item(97, 13).Value = -346.4737
$ws.Cells.Item(102, 8).Value = 387377.84
$ws.Cells.Item(102, 9).Value = 541354.4
$ws.Cells.Item(102, 11).Value = 541354.4
$ws.Cells.Item(102, 13).Value = -539732.4
$ws.Cells.Item(132, 8).Value = 2843.1843
$ws.Cells.Item(132, 9).Value = 2419.2942
$ws.Cells.Item(132, 11).Value = 7257.882599999999
$ws.Cells.Item(132, 13).Value = -4727.882599999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 11880.091
$ws.Cells.Item(20, 9).Value = 2987.1667
$ws.Cells.Item(20, 11).Value = 2987.1667
$ws.Cells.Item(20, 13).Value = -2740.1667
$ws.Cells.Item(86, 8).Value = 4397.2354
$ws.Cells.Item(86, 9).Value = 2305.3635
$ws.Cells.Item(86, 11).Value = 2305.3635
$ws.Cells.Item(86, 13).Value = -1182.3635
$ws.Cells.Item(89, 8).Value = 4397.2354
$ws.Cells.Item(89, 9).Value = 2305.3635
$ws.Cells.Item(89, 11).Value = 11526.8175
$ws.Cells.Item(89, 13).Value = -5910.817499999999
$ws.Cells.Item(94, 8).Value = 8621809
$ws.Cells.Item(94, 9).Value = 9260436
$ws.Cells.Item(94, 11).Value = 9260436
$ws.Cells.Item(94, 13).Value = -9259985
$ws.Cells.Item(99, 8).Value = 2851.55
$ws.Cells.Item(99, 9).Value = 2795.9412
$ws.Cells.Item(99, 11).Value = 2795.9412
$ws.Cells.Item(99, 13).Value = -1297.9412
$ws.Cells.Item(105, 8).Value = 2059.0667
$ws.Cells.Item(105, 9).Value = 1945.1538
$ws.Cells.Item(105, 10).Value = 2799.5
$ws.Cells.Item(105, 11).Value = 1945.1538
$ws.Cells.Item(105, 12).Value = 2799.5
$ws.Cells.Item(105, 13).Value = -198.1538
$ws.Cells.Item(105, 14).Value = -6293.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4010.6572
$ws.Cells.Item(31, 9).Value = 2911.4614
$ws.Cells.Item(31, 11).Value = 2911.4614
$ws.Cells.Item(31, 13).Value = -2616.4614
$ws.Cells.Item(34, 8).Value = 4010.6572
$ws.Cells.Item(34, 9).Value = 2911.4614
$ws.Cells.Item(34, 11).Value = 2911.4614
$ws.Cells.Item(34, 13).Value = -2709.4614
$ws.Cells.Item(107, 8).Value = 901.1111
$ws.Cells.Item(107, 9).Value = 951.65
$ws.Cells.Item(107, 10).Value = 756.7143
$ws.Cells.Item(107, 11).Value = 951.65
$ws.Cells.Item(107, 12).Value = 756.7143
$ws.Cells.Item(107, 13).Value = 968.35
$ws.Cells.Item(107, 14).Value = -4596.7143
$ws.Cells.Item(134, 8).Value = 3205.1206
$ws.Cells.Item(134, 9).Value = 2647.1135
$ws.Cells.Item(134, 11).Value = 7941.3405
$ws.Cells.Item(134, 13).Value = -5406.3405
$ws.Cells.Item(141, 8).Value = 552419.2
$ws.Cells.Item(141, 10).Value = 552419.2
$ws.Cells.Item(141, 12).Value = 552419.2
$ws.Cells.Item(141, 14).Value = -562779.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 518.7778
$ws.Cells.Item(7, 9).Value = 565
$ws.Cells.Item(7, 10).Value = 495.66666
$ws.Cells.Item(7, 11).Value = 1695
$ws.Cells.Item(7, 12).Value = 1486.99998
$ws.Cells.Item(7, 13).Value = -1583
$ws.Cells.Item(7, 14).Value = -1710.99998
$ws.Cells.Item(34, 8).Value = 2285.6428
$ws.Cells.Item(34, 9).Value = 1000
$ws.Cells.Item(34, 10).Value = 2384.5386
$ws.Cells.Item(34, 11).Value = 3000
$ws.Cells.Item(34, 12).Value = 7153.6158
$ws.Cells.Item(34, 13).Value = -2916
$ws.Cells.Item(34, 14).Value = -7321.6158
$ws.Cells.Item(63, 8).Value = 22442.723
$ws.Cells.Item(63, 9).Value = 23750
$ws.Cells.Item(63, 10).Value = 11984.5
$ws.Cells.Item(63, 11).Value = 71250
$ws.Cells.Item(63, 12).Value = 35953.5
$ws.Cells.Item(63, 13).Value = -70501
$ws.Cells.Item(63, 14).Value = -37451.5
$ws.Cells.Item(64, 8).Value = 5014
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 5014
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 15042
$ws.Cells.Item(64, 13).ClearContents()
$ws.Cells.Item(64, 14).Value = -15582
$ws.Cells.Item(66, 8).Value = 22442.723
$ws.Cells.Item(66, 9).Value = 23750
$ws.Cells.Item(66, 10).Value = 11984.5
$ws.Cells.Item(66, 11).Value = 213750
$ws.Cells.Item(66, 12).Value = 107860.5
$ws.Cells.Item(66, 13).Value = -210006
$ws.Cells.Item(66, 14).Value = -115348.5
$ws.Cells.Item(67, 8).Value = 5014
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 10).Value = 5014
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 12).Value = 15042
$ws.Cells.Item(67, 13).ClearContents()
$ws.Cells.Item(67, 14).Value = -16914
$ws.Cells.Item(134, 8).Value = 5938.1665
$ws.Cells.Item(134, 9).Value = 943
$ws.Cells.Item(134, 11).Value = 2829
$ws.Cells.Item(134, 13).Value = 2241
$ws.Cells.Item(141, 8).Value = 5198
$ws.Cells.Item(141, 9).Value = 5872.5
$ws.Cells.Item(141, 11).Value = 17617.5
$ws.Cells.Item(141, 13).Value = -12437.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(75, 8).Value = 40000
$ws.Cells.Item(75, 10).Value = 40000
$ws.Cells.Item(75, 12).Value = 40000
$ws.Cells.Item(75, 14).Value = -41748
$ws.Cells.Item(78, 8).Value = 40000
$ws.Cells.Item(78, 10).Value = 40000
$ws.Cells.Item(78, 12).Value = 120000
$ws.Cells.Item(78, 14).Value = -128736
$ws.Cells.Item(113, 8).Value = 3421.2104
$ws.Cells.Item(113, 9).Value = 2696.3845
$ws.Cells.Item(113, 11).Value = 2696.3845
$ws.Cells.Item(113, 13).Value = -526.3845000000001
$ws.Cells.Item(132, 8).Value = 3676.6738
$ws.Cells.Item(132, 9).Value = 3111.6758
$ws.Cells.Item(132, 10).Value = 5999.4443
$ws.Cells.Item(132, 11).Value = 9335.027399999999
$ws.Cells.Item(132, 12).Value = 17998.3329
$ws.Cells.Item(132, 13).Value = -6805.027399999999
$ws.Cells.Item(132, 14).Value = -23058.3329
$ws.Cells.Item(134, 8).Value = 41297
$ws.Cells.Item(134, 10).Value = 41297
$ws.Cells.Item(134, 12).Value = 123891
$ws.Cells.Item(134, 14).Value = -128961

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 574.8823
$ws.Cells.Item(16, 9).Value = 604.5625
$ws.Cells.Item(16, 11).Value = 604.5625
$ws.Cells.Item(16, 13).Value = -434.5625
$ws.Cells.Item(68, 8).Value = 4118.4287
$ws.Cells.Item(68, 9).Value = 3354.8333
$ws.Cells.Item(68, 11).Value = 3354.8333
$ws.Cells.Item(68, 13).Value = -2605.8333
$ws.Cells.Item(71, 8).Value = 4118.4287
$ws.Cells.Item(71, 9).Value = 3354.8333
$ws.Cells.Item(71, 11).Value = 16774.1665
$ws.Cells.Item(71, 13).Value = -13030.1665
$ws.Cells.Item(82, 8).Value = 2056
$ws.Cells.Item(82, 9).Value = 1971.625
$ws.Cells.Item(82, 11).Value = 1971.625
$ws.Cells.Item(82, 13).Value = -1610.625
$ws.Cells.Item(85, 8).Value = 2056
$ws.Cells.Item(85, 9).Value = 1971.625
$ws.Cells.Item(85, 11).Value = 1971.625
$ws.Cells.Item(85, 13).Value = -723.625
$ws.Cells.Item(132, 8).Value = 14497308
$ws.Cells.Item(132, 9).Value = 18521644
$ws.Cells.Item(132, 11).Value = 55564932
$ws.Cells.Item(132, 13).Value = -55562402

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 9618443
$ws.Cells.Item(100, 9).Value = 11908349
$ws.Cells.Item(100, 11).Value = 23816698
$ws.Cells.Item(100, 13).Value = -23816157
$ws.Cells.Item(107, 8).Value = 342.85715
$ws.Cells.Item(107, 9).Value = 342.5
$ws.Cells.Item(107, 11).Value = 1027.5
$ws.Cells.Item(107, 13).Value = 892.5
$ws.Cells.Item(132, 8).Value = 445458.47
$ws.Cells.Item(132, 9).Value = 694349.1
$ws.Cells.Item(132, 10).Value = 20880.234
$ws.Cells.Item(132, 11).Value = 2083047.3
$ws.Cells.Item(132, 12).Value = 62640.702
$ws.Cells.Item(132, 13).Value = -2080517.3
$ws.Cells.Item(132, 14).Value = -67700.702
$ws.Cells.Item(133, 8).Value = 27665.666
$ws.Cells.Item(133, 10).Value = 26123.875
$ws.Cells.Item(133, 12).Value = 26123.875
$ws.Cells.Item(133, 14).Value = -36243.875

Write-Output "Applied 237 cell updates and 2 clears across 8 sheets"